# Weekly update: insert a new weekly record at row 27 for
# "Vega Monumental Concepción" / Jengibre, pushing the existing
# historical rows (27-44) down by one (to 28-45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 27. This shifts the
# existing rows 27-44 down to 28-45 (values and styles included).
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Cells.Item(27, 1).Value2 = 11
$ws.Cells.Item(27, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(27, 3).Value2 = "Bíobío"
$ws.Cells.Item(27, 4).Value2 = 44978
$ws.Cells.Item(27, 5).Value2 = 8
$ws.Cells.Item(27, 6).Value2 = 100114007
$ws.Cells.Item(27, 7).Value2 = "Jengibre"
$ws.Cells.Item(27, 8).Value2 = "Sin especificar"
$ws.Cells.Item(27, 9).Value2 = "Primera"
$ws.Cells.Item(27, 10).Value2 = 40
$ws.Cells.Item(27, 11).Value2 = 13000
$ws.Cells.Item(27, 12).Value2 = 14000
$ws.Cells.Item(27, 13).Value2 = 13500
$ws.Cells.Item(27, 14).Value2 = "`$/caja 13 kilos"
$ws.Cells.Item(27, 15).Value2 = "Perú"
$ws.Cells.Item(27, 16).Value2 = 1038
$ws.Cells.Item(27, 17).Value2 = 13
$ws.Cells.Item(27, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same date number format used by
# the rest of column D.
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
